$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep A1 "touched" (format-only, no value/new style) so the sheet's used-range /
# dimension keeps starting at row 1, matching the original authoring tool's output.
$f = $ws.Range("A1").Font
$f.Name = $f.Name

# The workbook used to show 5 duplicate "employee" blocks across columns A-AS on
# row 2/3 (date header) and row 5 (shift header). We are trimming this down to
# just two blocks (A-I and J-R) and relabeling the shift names as "Presize"
# shifts instead of the old Cherry Line / Operations shift names.

# Before removing the 3rd/4th/5th blocks (S-AS), copy the formatting (fills,
# borders, fonts) from the last block (AB:AJ) onto the 2nd block (J:R) so the
# remaining 2nd/"last" block keeps the distinctive "last column" styling.
$ws.Range("AB2:AJ3").Copy()
$ws.Range("J2:R3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("AB5:AJ5").Copy()
$ws.Range("J5:R5").PasteSpecial(-4122)   # xlPasteFormats

# Relabel the shift names.
$ws.Range("A5").Value2 = "Presize 7:00AM - 3:30PM"
$ws.Range("J5").Value2 = "Presize 4:00PM - 12:30AM"

# Remove the now-unneeded employee columns S through AS (3 extra blocks).
$ws.Range("S1:AS1").EntireColumn.Delete()

# Restore the selection to the remaining last header cell.
$ws.Range("J5").Select()
